# Auto-generated edit script: apply cryptos price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.025.05"
$ws.Range("E2").Value = "  +1.72%  "

# Row 3
$ws.Range("D3").Value = "3.424.53"
$ws.Range("E3").Value = "  +1.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.69%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  +0.19%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.68%  "

# Row 10
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.08%  "

# Row 12
$ws.Range("D12").Value = "4.012.91"
$ws.Range("E12").Value = "  +1.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.28%  "

# Row 14
$ws.Range("E14").Value = "  -0.56%  "

# Row 15
$ws.Range("D15").Value = "3.422.03"
$ws.Range("E15").Value = "  +0.76%  "

# Row 16
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
$ws.Range("D17").Value = "62.156.25"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.47%  "

# Row 19
$ws.Range("E19").Value = "  +3.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.70%  "

# Row 22
$ws.Range("E22").Value = "  -1.81%  "

# Row 23
$ws.Range("E23").Value = "  +0.80%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("E25").Value = "  +0.81%  "

# Row 26
$ws.Range("D26").Value = "3.567.75"
$ws.Range("E26").Value = "  +1.33%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.49%  "

# Row 28
$ws.Range("E28").Value = "  +4.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

# Row 31
$ws.Range("E31").Value = "  +0.80%  "

# Row 32
$ws.Range("E32").Value = "  +3.02%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.21%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.72%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.42%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.88%  "

# Row 39
$ws.Range("D39").Value = "3.458.83"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "29.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0755"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.788"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.44%  "

# Row 43
$ws.Range("E43").Value = "  +2.31%  "

# Row 44
$ws.Range("E44").Value = "  +1.99%  "

# Row 45
$ws.Range("E45").Value = "  +4.64%  "

# Row 46
$ws.Range("D46").Value = "2.512.29"
$ws.Range("E46").Value = "  +2.22%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.08%  "

# Row 49
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0265"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "

# Row 51
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "

